$d = $word.ActiveDocument

# 1. Remove the "Data do relatório: ..." paragraph together with the
#    single-space paragraph that immediately follows it. Locate the
#    paragraph by its text so the edit is anchored to content rather
#    than a hard-coded paragraph index.
$dateRng = $d.Content
$dateRng.Find.Execute("Data do relatório: 22 de janeiro de 2024", $true, $false, $false,
                       $false, $false, $true, 1, $false, "", 0)
$startPara = $dateRng.Paragraphs(1)
$endPara = $startPara.Next()
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()

# 2. Reword the closing sentence of the document.
#    The run being edited sits right after another run that happens to
#    share identical run formatting (a lone space). A plain Find/Replace
#    (or a direct Range.Text assignment) causes the two runs to coalesce
#    into one, which would also swallow that separator space into the
#    new text. Nudging the run's formatting away from its neighbour
#    before the edit - then restoring it on the freshly written text -
#    keeps the edited run distinct, matching the original run layout.
$oldSentence = "Se o produto conseguir manter seu ritmo atual de vendas no epicentro do mundo da saúde e boa forma, então ele poderá estar pronto para uma exposição nacional."
$newSentence = "Se o produto puder manter seu ritmo atual de vendas no mundo da saúde e fitness, ele poderá estar pronto para exposição nacional."

$found = $d.Content
$found.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$origColor = $found.Font.Color
$found.Font.Color = 255
$found.Text = $newSentence

$replaced = $d.Content
$replaced.Find.Execute($newSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$replaced.Font.Color = $origColor
